# Descripcion-DiagramadeClase.docx correction:
# the class "Consumible" was renamed to "Fruta" in the elaboration text
# (bullets 2 and 3 of the bulleted list).

$d = $word.ActiveDocument

function Replace-WordAsNewRun($findText, $newText) {
    # Locate the text to replace.
    $found = $d.Content
    $found.Find.Execute($findText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)
    $start = $found.Start
    $end = $found.End

    # Remove the old text and insert the new text in its place.
    $d.Range($start, $end).Text = ""
    $d.Range($start, $start).InsertAfter($newText)

    # Nudge formatting on the freshly-inserted text so it becomes its own
    # run instead of being silently re-absorbed into the neighbouring run.
    $newRun = $d.Range($start, $start + $newText.Length)
    $newRun.Bold = 1
    $newRun.Bold = 0
}

# Bullet 2: "Las clases Consumible y Pokebola heredan de Bonus, ..."
#        -> "Las clases Fruta y Pokebola heredan de Bonus, ..."
Replace-WordAsNewRun "Consumible" "Fruta"

# Bullet 3: "La clase costanera instancia Personaje, Pokebola y Consumible. "
#        -> "La clase costanera instancia Personaje, Pokebola y Fruta."
Replace-WordAsNewRun "Consumible. " "Fruta."
